$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.347.17"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.16"
$ws.Range("E3").Value = "  +0.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.56"
$ws.Range("E5").Value = "  +0.56%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.575"
$ws.Range("E6").Value = "  +3.67%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.10"
$ws.Range("E8").Value = "  +10.50%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.301"
$ws.Range("E9").Value = "  +2.23%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("E10").Value = "  +0.69%  "

# Row 11
$ws.Range("E11").Value = "  +1.97%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.064.38"
$ws.Range("E12").Value = "  +0.85%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.75"
$ws.Range("E13").Value = "  +6.58%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.876.24"
$ws.Range("E14").Value = "  +4.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.644"
$ws.Range("E15").Value = "  +1.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.50"
$ws.Range("E16").Value = "  +5.00%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.300.44"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.03"
$ws.Range("E18").Value = "  +1.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.70"
$ws.Range("E19").Value = "  +0.44%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0794"
$ws.Range("E20").Value = "  +0.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.62"
$ws.Range("E21").Value = "  +3.86%  "

# Row 22
$ws.Range("E22").Value = "  +0.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("E23").Value = "  +0.75%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "171.85"
$ws.Range("E24").Value = "  +2.92%  "

# Row 25
$ws.Range("E25").Value = "  +2.96%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.97"
$ws.Range("E26").Value = "  +9.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.89"
$ws.Range("E27").Value = "  +2.19%  "

# Row 28
$ws.Range("E28").Value = "  +2.58%  "

# Row 29
$ws.Range("E29").Value = "  -0.11%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.02"
$ws.Range("E30").Value = "  +0.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0532"
$ws.Range("E31").Value = "  +1.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.85"
$ws.Range("E32").Value = "  +0.98%  "

# Row 33
$ws.Range("E33").Value = "  +0.58%  "

# Row 34
$ws.Range("E34").Value = "  +0.61%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.390.22"
$ws.Range("E35").Value = "  -1.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.673"
$ws.Range("E36").Value = "  -1.35%  "

# Row 37
$ws.Range("E37").Value = "  -5.79%  "

# Row 38
$ws.Range("E38").Value = "  -0.43%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0191"
$ws.Range("E39").Value = "  +0.02%  "

# Row 40
$ws.Range("E40").Value = "  +10.74%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.963"
$ws.Range("E41").Value = "  +2.68%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("E44").Value = "  +0.40%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.45"
$ws.Range("E45").Value = "  -2.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.01"
$ws.Range("E46").Value = "  -0.83%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0503"
$ws.Range("E47").Value = "  -4.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.964.01"
$ws.Range("E48").Value = "  +0.83%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.67"
$ws.Range("E49").Value = "  -0.65%  "

# Row 50
$ws.Range("E50").Value = "  +0.04%  "

# Row 51
$ws.Range("E51").Value = "  -0.38%  "

# Row 42 - was MXToken, now Aave
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "82.11"
$ws.Range("E42").Value = "  -2.27%  "

# Row 43 - was Aave, now MXToken
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.82"
$ws.Range("E43").Value = "  +1.08%  "
